$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# New handback record identifiers (this run)
# ----------------------------------------------------------------------
$uuid1 = "77879806-9de0-495c-b6c6-b7169e833960"
$uuid2 = "87db403b-203f-44e4-a0e0-0276bc3326ca"
$hash1 = "bee9ce64a15931c20b8ced5c1de9ba1434e61f1a"
$hash2 = "e648f48d8489fd8853d6d56454fe0e325d8811ee"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Delete()

# Row 2 -- existing handback record, refreshed with this run's data
$ov.Range("A2").Value = "$uuid1.md"
$ov.Range("B2").Value = "e2e\$uuid1.md"
$ov.Range("C2").Value = ".md"
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("G2").NumberFormat = $dateFmt
$ov.Range("G2").Value = "2016-09-06 23:16:00"

# Row 3 -- new handback record
$ov.Range("A3").Value = "$uuid2.md"
$ov.Range("B3").Value = "e2e\$uuid2.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").NumberFormat = $dateFmt
$ov.Range("G3").Value = "2016-09-06 23:16:00"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/beaa13ef3e5483f1c1da9f5f50e4d513f01776b9/e2e/$uuid1.md", "", "", "e2e\$uuid1.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/beaa13ef3e5483f1c1da9f5f50e4d513f01776b9/e2e/$uuid2.md", "", "", "e2e\$uuid2.md") | Out-Null

$ovList = $ov.ListObjects.Item(1)
$ovList.Resize($ov.Range("A1:G3"))

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

# Row 2 -- existing handback record, refreshed with this run's data
$zh.Range("A2").Value = "$uuid1.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
$zh.Range("F2").Value = "'False"
$zh.Range("G2").Value = "$uuid1.$hash1.zh-cn.xlf"
$zh.Range("H2").NumberFormat = $dateFmt
$zh.Range("H2").Value = "2016-09-06 23:15:54"
$zh.Range("I2").Value = "$uuid1.md"
$zh.Range("J2").Value = "$uuid1.$hash1.zh-cn.xlf"
$zh.Range("K2").NumberFormat = $dateFmt
$zh.Range("K2").Value = "2016-09-06 23:16:31"
$zh.Range("L2").Value = "'"
$zh.Range("M2").Value = "'True"
$zh.Range("N2").Value = "'"
$zh.Range("O2").Value = "'False"
$zh.Range("P2").Value = "'"

# Row 3 -- new handback record
$zh.Range("A3").Value = "$uuid2.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "$uuid2.$hash2.zh-cn.xlf"
$zh.Range("H3").NumberFormat = $dateFmt
$zh.Range("H3").Value = "2016-09-06 23:15:54"
$zh.Range("I3").Value = "$uuid2.md"
$zh.Range("J3").Value = "$uuid2.$hash2.zh-cn.xlf"
$zh.Range("K3").NumberFormat = $dateFmt
$zh.Range("K3").Value = "2016-09-06 23:16:31"
$zh.Range("L3").Value = "'"
$zh.Range("M3").Value = "'True"
$zh.Range("N3").Value = "'"
$zh.Range("O3").Value = "'False"
$zh.Range("P3").Value = "'"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5bd4fa5d42dc047217cdb38d021b220ee3509b2f/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5bd4fa5d42dc047217cdb38d021b220ee3509b2f/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5bd4fa5d42dc047217cdb38d021b220ee3509b2f/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5bd4fa5d42dc047217cdb38d021b220ee3509b2f/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null

$zhList = $zh.ListObjects.Item(1)
$zhList.Resize($zh.Range("A1:P3"))

# =======================================================================
# Sheet "de-de"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

# Row 2 -- existing handback record, refreshed with this run's data
$de.Range("A2").Value = "$uuid1.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
$de.Range("F2").Value = "'False"
$de.Range("G2").Value = "$uuid1.$hash1.de-de.xlf"
$de.Range("H2").NumberFormat = $dateFmt
$de.Range("H2").Value = "2016-09-06 23:16:00"
$de.Range("I2").Value = "$uuid1.md"
$de.Range("J2").Value = "$uuid1.$hash1.de-de.xlf"
$de.Range("K2").NumberFormat = $dateFmt
$de.Range("K2").Value = "2016-09-06 23:16:40"
$de.Range("L2").Value = "'"
$de.Range("M2").Value = "'True"
$de.Range("N2").Value = "'"
$de.Range("O2").Value = "'False"
$de.Range("P2").Value = "'"

# Row 3 -- new handback record
$de.Range("A3").Value = "$uuid2.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "$uuid2.$hash2.de-de.xlf"
$de.Range("H3").NumberFormat = $dateFmt
$de.Range("H3").Value = "2016-09-06 23:16:00"
$de.Range("I3").Value = "$uuid2.md"
$de.Range("J3").Value = "$uuid2.$hash2.de-de.xlf"
$de.Range("K3").NumberFormat = $dateFmt
$de.Range("K3").Value = "2016-09-06 23:16:40"
$de.Range("L3").Value = "'"
$de.Range("M3").Value = "'True"
$de.Range("N3").Value = "'"
$de.Range("O3").Value = "'False"
$de.Range("P3").Value = "'"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5d86c42e5dc66d80d862d18f4557db40b716142/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5d86c42e5dc66d80d862d18f4557db40b716142/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5d86c42e5dc66d80d862d18f4557db40b716142/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5d86c42e5dc66d80d862d18f4557db40b716142/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null

$deList = $de.ListObjects.Item(1)
$deList.Resize($de.Range("A1:P3"))

# =======================================================================
# Apply the same font styling used for hyperlink cells (A/I columns) and
# keep the original column-header cell look for row 1 untouched.
# =======================================================================
foreach ($sheetInfo in @(
        @{ Sheet = $ov; Cells = @("B3") },
        @{ Sheet = $zh; Cells = @("A3", "I3") },
        @{ Sheet = $de; Cells = @("A3", "I3") }
    )) {
    foreach ($addr in $sheetInfo.Cells) {
        $rng = $sheetInfo.Sheet.Range($addr)
        $rng.Font.Underline = $true
        $rng.Font.Color = 15570276
    }
}
